$wb = $excel.ActiveWorkbook

# Duplicate the "Sheet5" (MNIST) analysis sheet to create a new "Sheet6"
# (MNIST Kaggle) sheet, placed after all existing sheets.
$src = $wb.Worksheets.Item("Sheet5")
$src.Activate()
$src.Range("A1:H13").Select()
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "Sheet6"

# Re-title the copy and drop the trailing spacer row that belonged to Sheet5.
$new.Range("A1").Value = "MNIST Kaggle"
$new.Rows("14:14").Delete()

$new.Activate()
$new.Range("A8").Select()
